# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-07 06:22:31
#
# The "Recorded By" column (G) lists who recorded/edited a session, e.g.
#   "System, dnasr281@gmail.com"
# A few rows had the "System" entry re-ordered relative to the other
# recorder(s) in the comma separated list. Re-apply that re-ordering here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($text -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
